# Remove invalid benchmarking rows where the hardware configuration violates
# n_nodes validation rules (n_mpi_procs must be >= n_nodes; a single task
# cannot span multiple nodes). These rows (by their original 1-based sheet
# row numbers) are no longer valid combinations and are deleted, with the
# remaining rows shifting up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = 50, 49, 48, 30, 28, 23
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
